# Update countries & provincias Spain
# - Refresh COVID numbers for a handful of countries (incl. swapping the
#   Montserrat / Islas Malvinas rows so Islas Malvinas sorts above Montserrat)
# - Bump the "Datos actualizados" timestamp from 13:27 to 14:44

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp banner in row 1
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 14:44"

# row -> [Pais, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes]
$updates = @{
    4   = @("Estados Unidos", 5656744, 770,  3012050, 2469589, 0, 31,  175105)
    14  = @("Iran",            350279, 2444, 302528,  27626,   0, 153, 20125)
    36  = @("Suecia",          85411,  0,    0,       0,       0, 5,   5802)
    41  = @("Kuwait",          78145,  675,  69771,   7867,    0, 2,   507)
    45  = @("Paises Bajos",    64525,  552,  0,       0,       0, 6,   6181)
    61  = @("Uzbekistan",      36968,  616,  32557,   4164,    0, 5,   247)
    80  = @("Dinamarca",       15940,  85,   13651,   1668,    0, 0,   621)
    102 = @("Croacia",         7074,   219,  5386,    1520,    0, 2,   168)
    213 = @("Islas Malvinas",  13,     0,    13,      0,       0, 0,   0)
    214 = @("Montserrat",      13,     0,    12,      0,       0, 0,   1)
}

foreach ($r in $updates.Keys) {
    $row = $updates[$r]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $ws.Range("H$r").Value = $row[7]
}
